$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, establish row 53 with the correct style by copying the date-column
# format from row 52 (column A) before populating values.
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A53").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

# Update data rows 2-53 with recomputed forecast values. Row 2 introduces a new
# leading 2007 data point, shifting the remaining historical forecasts down by
# one row relative to the previous version, and row 53 is a brand-new trailing
# data point for 2025/2026.
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = -0.7015558851707349
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 0.3143490788445336
$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 0.7771450785698075
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 0.4748521911469794
$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = 0.3590181115727287
$ws.Range("D4").Value = 2009
$ws.Range("E4").Value = 0.6970543652217165
$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = 2009
$ws.Range("C5").Value = 0.9995687521967556
$ws.Range("D5").Value = 2010
$ws.Range("E5").Value = 0.02446583860156171
$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = -0.01587181126743165
$ws.Range("D6").Value = 2010
$ws.Range("E6").Value = -0.1225239755399454
$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = 0.5978820435291077
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = 0.7697980859487474
$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = -0.02256889165885845
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = -0.02753509623226735
$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = -0.5061359875450311
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = -0.362324052998142
$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 0.09611428386597787
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = -0.2256894584805158
$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = -0.2706540469742502
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = -0.06282556559906727
$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = -0.1827723404408288
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = -0.2251688766575
$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = 0.006024133679316535
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 0.01247916696665019
$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = -0.001350220946483294
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = 0.1494097328869959
$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = 0.3981709080043139
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 0.1740860482467133
$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = -0.075754880139145
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = -0.5497151367044428
$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = -0.2945738319855118
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = -0.03768624985649449
$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = -0.5761528471665445
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = -0.1255150964614482
$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = -0.03047919532177534
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = -0.1249617237519152
$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = -0.2011999787958185
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = -0.3746351385105373
$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 0.1932702877606385
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = -0.250093582508859
$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 0.1213692818849532
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 0.1731436979489631
$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 0.2738179272064434
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 0.2001500500062425
$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = -0.1256133802673975
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = -0.02501876407304815
$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 0.08160407806372394
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 0.2722769713627837
$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 0.1493219406571766
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = 0.9013851022877439
$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 0.2799548089016612
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 0.3604862916655627
$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = -0.2736870064301455
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = -0.11321783823105
$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = -0.2916219766884276
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = -0.2643564666883758
$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = -0.4278219446121612
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = -0.8258413506386342
$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = -0.3225448117294083
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = -0.2397840863870626
$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = -0.2932081122163033
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = -0.1126446518617819
$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = -0.4412356890029168
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = -0.3695048299469872
$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = -1.026566979837418
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = -2.816143384276215
$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = -0.7059330669460406
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = -1.039648953489214
$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = -0.03584227163500042
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = -0.1719585843969917
$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = -0.04976849661378902
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = -1.608095448141689
$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = 0.3179894933462268
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = 0.07011423530434158
$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = 0.5748854066360609
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = 0.35467759793264
$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 1.323454226677478
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 0.7858329241748896
$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = 0.4784173072842179
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 0.1043052025668345
$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = 0.463604920919658
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = -0.3203420516749933
$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = 0.6932214722757513
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 0.9777431256012825
$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = 0.8831924739260089
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 0.8355283619122744
$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = 0.4279153732809959
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = 0.02757741937535751
$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = 0.621639092134818
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = 0.009137938461889483
$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = -0.3894584472036278
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = 0.1876441418131369
$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = -1.081515348061801
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = -0.07003400812273242
$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = -0.9378224616154895
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = 0.3237070386496343
$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = -0.6768900623516982
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = 0.9453792747973422
$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = 1.909988706581967
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = 0.5955700148392751
$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = 0.4626514211933497
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = -0.4747835872719319
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = 1.136769786738334
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = 0.294850926654866
